# SIGCOMM poster.pptx — "Don't add newline to graphviz format"
#
# The bullet describing the Propane compilation stage was retyped in
# place. The resulting text is unchanged, but PowerPoint now stores it
# as four runs (split at the word boundaries that were touched) instead
# of one. Reproduce that by re-keying the two interior chunks, which
# forces PowerPoint to split the original single run into four runs
# that carry the same character formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate "Text Box 194" by name (its index has been 9 historically, but
# search by name first so this keeps working if shapes get reordered).
$shp = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "Text Box 194") {
        $shp = $candidate
        break
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(9)
}

$tr = $shp.TextFrame.TextRange

# Find the paragraph that holds the bullet text (search, in case the
# index ever shifts) instead of hard-coding paragraph 3.
$target = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text.StartsWith("Compile Propane to per-destination")) {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the 'Compile Propane...' bullet paragraph"
}

# The textbox has AutoFit-shape-to-text turned on, so rewriting its
# runs can make the host recompute the shape height. Remember it and
# put it back once the text is settled so geometry stays untouched,
# matching the diff (which only touches <a:t> content).
$origHeight = $shp.Height

# Re-key the two interior word groups in place; each re-assignment
# splits a run boundary at that point without touching the rest of the
# text, ending up with four runs whose concatenated text equals the
# original sentence. Characters() on a paragraph-scoped TextRange is
# 1-based and relative to the paragraph itself (not the text frame).
$chunk1 = "Compile Propane to per-destination "
$chunk2 = "state "
$chunk3 = "machines that "
$chunk4 = "associate paths with ranks. A lower rank means the path is preferred."

$offset2 = $chunk1.Length + 1
$offset3 = $offset2 + $chunk2.Length

$target.Characters($offset2, $chunk2.Length).Text = $chunk2
$target.Characters($offset3, $chunk3.Length).Text = $chunk3

$shp.Height = $origHeight

# Paragraphs()'s .Text includes the trailing paragraph-mark (CR), so
# trim that before comparing against the plain concatenated chunks.
$expected = $chunk1 + $chunk2 + $chunk3 + $chunk4
$actual = $target.Text.TrimEnd("`r")
if ($actual -ne $expected) {
    throw "Unexpected paragraph text after edit: $actual"
}

Write-Host "Final paragraph text:" $target.Text
